# Scheduled-runner data refresh: update market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, cols H-N) on
# affected rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H62").Value = 6284.9546
$ws_ALC.Range("I62").Value = 6237.2
$ws_ALC.Range("J62").Value = 6387.2856
$ws_ALC.Range("K62").Value = 6237.2
$ws_ALC.Range("L62").Value = 6387.2856
$ws_ALC.Range("M62").Value = -5613.2
$ws_ALC.Range("N62").Value = -7635.2856
$ws_ALC.Range("H65").Value = 6284.9546
$ws_ALC.Range("I65").Value = 6237.2
$ws_ALC.Range("J65").Value = 6387.2856
$ws_ALC.Range("K65").Value = 31186
$ws_ALC.Range("L65").Value = 31936.428
$ws_ALC.Range("M65").Value = -28066
$ws_ALC.Range("N65").Value = -38176.428
$ws_ALC.Range("H76").Value = 2500
$ws_ALC.Range("I76").Value = 2500
$ws_ALC.Range("K76").Value = 2500
$ws_ALC.Range("M76").Value = -2185
$ws_ALC.Range("H79").Value = 2500
$ws_ALC.Range("I79").Value = 2500
$ws_ALC.Range("K79").Value = 2500
$ws_ALC.Range("M79").Value = -1408
$ws_ALC.Range("H106").Value = 4484.75
$ws_ALC.Range("I106").Value = 1939
$ws_ALC.Range("J106").Value = 5333.3335
$ws_ALC.Range("K106").Value = 1939
$ws_ALC.Range("L106").Value = 5333.3335
$ws_ALC.Range("M106").Value = -1308
$ws_ALC.Range("N106").Value = -6595.3335
$ws_ALC.Range("H132").Value = 27786456
$ws_ALC.Range("I132").Value = 34484868
$ws_ALC.Range("K132").Value = 103454604
$ws_ALC.Range("M132").Value = -103452074
$ws_ALC.Range("H137").Value = 3736.1667
$ws_ALC.Range("I137").Value = 3943.8
$ws_ALC.Range("K137").Value = 11831.4
$ws_ALC.Range("M137").Value = -9281.400000000001
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 7001.1963
$ws_ARM.Range("I32").Value = 7126.396
$ws_ARM.Range("K32").Value = 7126.396
$ws_ARM.Range("M32").Value = -6839.396
$ws_ARM.Range("H61").Value = 6313.372
$ws_ARM.Range("I61").Value = 4041.3794
$ws_ARM.Range("K61").Value = 4041.3794
$ws_ARM.Range("M61").Value = -3829.3794
$ws_ARM.Range("H74").Value = 2829.15
$ws_ARM.Range("I74").Value = 2304.861
$ws_ARM.Range("K74").Value = 2304.861
$ws_ARM.Range("M74").Value = -1430.861
$ws_ARM.Range("H77").Value = 2829.15
$ws_ARM.Range("I77").Value = 2304.861
$ws_ARM.Range("K77").Value = 11524.305
$ws_ARM.Range("M77").Value = -7156.305
$ws_ARM.Range("H136").Value = 6313.372
$ws_ARM.Range("I136").Value = 4041.3794
$ws_ARM.Range("K136").Value = 12124.1382
$ws_ARM.Range("M136").Value = -9574.138199999999
$ws_ARM.Range("H138").Value = 84499.5
$ws_ARM.Range("J138").Value = 84499.5
$ws_ARM.Range("L138").Value = 84499.5
$ws_ARM.Range("N138").Value = -94779.5
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H69").Value = 0
$ws_BSM.Range("J69").Value = 0
$ws_BSM.Range("L69").Value = 0
$ws_BSM.Range("N69").Value = $null
$ws_BSM.Range("H72").Value = 0
$ws_BSM.Range("J72").Value = 0
$ws_BSM.Range("L72").Value = 0
$ws_BSM.Range("N72").Value = $null
$ws_BSM.Range("H105").Value = 4491.4736
$ws_BSM.Range("I105").Value = 2967
$ws_BSM.Range("K105").Value = 2967
$ws_BSM.Range("M105").Value = -1220
$ws_BSM.Range("H107").Value = 2099
$ws_BSM.Range("I107").Value = 1723.75
$ws_BSM.Range("J107").Value = 2474.25
$ws_BSM.Range("K107").Value = 1723.75
$ws_BSM.Range("L107").Value = 2474.25
$ws_BSM.Range("M107").Value = 196.25
$ws_BSM.Range("N107").Value = -6314.25
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 2699.1086
$ws_CRP.Range("I31").Value = 1998.68
$ws_CRP.Range("K31").Value = 1998.68
$ws_CRP.Range("M31").Value = -1703.68
$ws_CRP.Range("H34").Value = 2699.1086
$ws_CRP.Range("I34").Value = 1998.68
$ws_CRP.Range("K34").Value = 1998.68
$ws_CRP.Range("M34").Value = -1796.68
$ws_CRP.Range("H105").Value = 367.5
$ws_CRP.Range("I105").Value = 367.5
$ws_CRP.Range("K105").Value = 367.5
$ws_CRP.Range("M105").Value = 1379.5
$ws_CRP.Range("H127").Value = 40416.625
$ws_CRP.Range("J127").Value = 40416.625
$ws_CRP.Range("L127").Value = 40416.625
$ws_CRP.Range("N127").Value = -50336.625
$ws_CRP.Range("H134").Value = 3122.0894
$ws_CRP.Range("I134").Value = 1277.8975
$ws_CRP.Range("K134").Value = 3833.6925
$ws_CRP.Range("M134").Value = -1298.6925
$ws_CRP.Range("H141").Value = 265405.8
$ws_CRP.Range("J141").Value = 265405.8
$ws_CRP.Range("L141").Value = 265405.8
$ws_CRP.Range("N141").Value = -275765.8
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 8105237
$ws_CUL.Range("J4").Value = 414812.66
$ws_CUL.Range("L4").Value = 1244437.98
$ws_CUL.Range("N4").Value = -1244661.98
$ws_CUL.Range("H5").Value = 9944.182000000001
$ws_CUL.Range("I5").Value = 548.8570999999999
$ws_CUL.Range("K5").Value = 1646.5713
$ws_CUL.Range("M5").Value = -1534.5713
$ws_CUL.Range("H131").Value = 1758.0857
$ws_CUL.Range("I131").Value = 842.36365
$ws_CUL.Range("K131").Value = 2527.09095
$ws_CUL.Range("M131").Value = 2512.90905
$ws_CUL.Range("H135").Value = 9944.182000000001
$ws_CUL.Range("I135").Value = 548.8570999999999
$ws_CUL.Range("K135").Value = 4939.7139
$ws_CUL.Range("M135").Value = -2404.7139
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H49").Value = 32792.5
$ws_GSM.Range("J49").Value = 32792.5
$ws_GSM.Range("L49").Value = 32792.5
$ws_GSM.Range("N49").Value = -33160.5
$ws_GSM.Range("H70").Value = 6853.1333
$ws_GSM.Range("I70").Value = 5724.625
$ws_GSM.Range("J70").Value = 8142.857
$ws_GSM.Range("K70").Value = 5724.625
$ws_GSM.Range("L70").Value = 8142.857
$ws_GSM.Range("M70").Value = -5454.625
$ws_GSM.Range("N70").Value = -8682.857
$ws_GSM.Range("H73").Value = 6853.1333
$ws_GSM.Range("I73").Value = 5724.625
$ws_GSM.Range("J73").Value = 8142.857
$ws_GSM.Range("K73").Value = 5724.625
$ws_GSM.Range("L73").Value = 8142.857
$ws_GSM.Range("M73").Value = -4788.625
$ws_GSM.Range("N73").Value = -10014.857
$ws_GSM.Range("H74").Value = 33333
$ws_GSM.Range("J74").Value = 33333
$ws_GSM.Range("L74").Value = 33333
$ws_GSM.Range("N74").Value = -35205
$ws_GSM.Range("H77").Value = 33333
$ws_GSM.Range("J77").Value = 33333
$ws_GSM.Range("L77").Value = 99999
$ws_GSM.Range("N77").Value = -109359
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H12").Value = 1500
$ws_LTW.Range("J12").Value = 1500
$ws_LTW.Range("L12").Value = 1500
$ws_LTW.Range("N12").Value = -1840
$ws_LTW.Range("H22").Value = 3443.7693
$ws_LTW.Range("J22").Value = 3982.4
$ws_LTW.Range("L22").Value = 3982.4
$ws_LTW.Range("N22").Value = -4572.4
$ws_LTW.Range("H25").Value = 31997.5
$ws_LTW.Range("I25").Value = 38995
$ws_LTW.Range("J25").Value = 25000
$ws_LTW.Range("K25").Value = 38995
$ws_LTW.Range("L25").Value = 25000
$ws_LTW.Range("M25").Value = -38765
$ws_LTW.Range("N25").Value = -25460
$ws_LTW.Range("H27").Value = 3443.7693
$ws_LTW.Range("J27").Value = 3982.4
$ws_LTW.Range("L27").Value = 3982.4
$ws_LTW.Range("N27").Value = -4196.4
$ws_LTW.Range("H95").Value = 47408
$ws_LTW.Range("J95").Value = 47408
$ws_LTW.Range("L95").Value = 47408
$ws_LTW.Range("N95").Value = -52900
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H101").Value = 62542.332
$ws_WVR.Range("J101").Value = 62542.332
$ws_WVR.Range("L101").Value = 62542.332
$ws_WVR.Range("N101").Value = -69032.33199999999
